# "Fix builk insert from input files"
# Bulk-insert two constant "language" columns (EN/English-Version and
# ES/Version-española) ahead of the existing per-row data, plus a new
# "line N" label column, pushing the original two data columns (a/b,
# c/d, e/f, g/h) out to columns F and G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original column A/B values for each row, now destined for F/G.
$pairs = @(
    @("a", "b"),
    @("c", "d"),
    @("e", "f"),
    @("g", "h")
)

for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "EN"
    $ws.Cells.Item($row, 2).Value = "English-Version"
    $ws.Cells.Item($row, 3).Value = "ES"
    $ws.Cells.Item($row, 4).Value = "Version-española"
    $ws.Cells.Item($row, 5).Value = "line $row"
    $ws.Cells.Item($row, 6).Value = $pairs[$i][0]
    $ws.Cells.Item($row, 7).Value = $pairs[$i][1]
}

# Column widths (character units), closest achievable values to the
# target worksheet's A:4.07 B:14.08 C:4.07 D:15.34 E:5.88 F:G 2.54.
$ws.Columns.Item(1).ColumnWidth = 3.17
$ws.Columns.Item(2).ColumnWidth = 13.17
$ws.Columns.Item(3).ColumnWidth = 3.17
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 5.0
$ws.Columns.Item(6).ColumnWidth = 1.67
$ws.Columns.Item(7).ColumnWidth = 1.67

$ws.Range("E4").Select()
